$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the C21 formula: remove the "+ 0.4 * C15/0.6" term
$ws.Range("C21").Formula = "=(C4-C19)*0.275"

# Update the active selection on the sheet view
$ws.Range("F26").Select()
